# Weekly price update: insert one new "Pepino ensalada" record for
# Comercializadora del Agro de Limarí as row 97, shifting the existing
# rows 97:182 down to 98:183 (dimension grows from R182 to R183).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 97 - pushes old rows 97..182 down to 98..183
# and extends the used range automatically.
$ws.Rows("97:97").Insert()

# Populate the newly inserted row with this week's data.
$ws.Range("A97").Value = 2
$ws.Range("B97").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C97").Value = "Coquimbo"
$ws.Range("D97").Value = 44902
$ws.Range("E97").Value = 4
$ws.Range("F97").Value = 100112043
$ws.Range("G97").Value = "Pepino ensalada"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 500
$ws.Range("K97").Value = 12000
$ws.Range("L97").Value = 13000
$ws.Range("M97").Value = 12500
$ws.Range("N97").Value = "$/caja 70 unidades"
$ws.Range("O97").Value = "Provincia de Limarí"
$ws.Range("P97").Value = 179
$ws.Range("Q97").Value = 70
$ws.Range("R97").Value = "Hortaliza"
